# edit.ps1 - Applies the two changes described by the diff:
#   1. Add a new "VenoBox" bullet item (ListParagraph, same numbering as the
#      other "Third Party Libraries Used" items) right after the "JQuery"
#      item, moving the trailing _GoBack bookmark onto the new paragraph.
#   2. Remove the stray spell-check proofErr wrapper (spellStart/spellEnd)
#      around "Timbo" in the table.

$d = $word.ActiveDocument

$wOpenXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wOpenXmlNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData>' +
           '</pkg:part></pkg:package>'
}

# --- 1. Insert "VenoBox" bullet after "JQuery" -----------------------------

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^JQuery\r$") {

        $jqueryPara = '<w:p w14:paraId="0D32513E" w14:textId="298DB3D5" w:rsidR="00FE6B7E" w:rsidRPr="00FE6B7E" w:rsidRDefault="00FE6B7E" w:rsidP="00FE6B7E">' +
                        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
                        '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>JQuery</w:t></w:r>' +
                      '</w:p>'

        $venoBoxPara = '<w:p>' +
                          '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
                          '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>VenoBox</w:t></w:r>' +
                          '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
                        '</w:p>'

        $target = $d.Range($p.Range.Start, $p.Range.End)
        $target.InsertXML((New-PkgXml ($jqueryPara + $venoBoxPara)))
        break
    }
}

# --- 2. Drop the proofErr spellStart/spellEnd wrapper around "Timbo" ------

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Timbo\r$") {

        $timboPara = '<w:p w14:paraId="6A05DF65" w14:textId="77777777" w:rsidR="00FD021E" w:rsidRDefault="00FD021E">' +
                       '<w:pPr><w:spacing w:line="240" w:lineRule="auto"/></w:pPr>' +
                       '<w:r><w:t>Timbo</w:t></w:r>' +
                     '</w:p>'

        $target = $d.Range($p.Range.Start, $p.Range.End)
        $target.InsertXML((New-PkgXml $timboPara))
        break
    }
}

Write-Host "Edits applied."
